$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# stay text (matching the source data, which stores these as text), without
# leaving a lingering number-format style on the cell.

$ws.Range("D2").Value = "57.217.20"
$ws.Range("E2").Value = "  +4.65%  "
$ws.Range("D3").Value = "2.360.17"
$ws.Range("E3").Value = "  +3.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("E8").Value = "  +1.67%  "
$ws.Range("D9").Value = "2.358.50"
$ws.Range("E9").Value = "  +2.13%  "
$ws.Range("E10").Value = "  +6.90%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.05%  "
$ws.Range("E13").Value = "  -0.14%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.779.78"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "57.059.46"
$ws.Range("E16").Value = "  +4.32%  "
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "2.367.24"
$ws.Range("E18").Value = "  +4.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.48%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("B25").Value = "Binance-PegBSC-USD"
$ws.Range("C25").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.161"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "0.0₃0746"
$ws.Range("E30").Value = "  +5.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.25"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("D37").Style = "Normal"
$ws.Range("E38").Value = "  +5.13%  "
$ws.Range("E39").Value = "  +8.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.18%  "
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("E42").Value = "  +5.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("E44").Value = "  +10.63%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.24%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("E49").Value = "  +5.35%  "
$ws.Range("E50").Value = "  +0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.80%  "
